$d = $word.ActiveDocument

# 1) "Leadership Style to be applied: [Insert your leadership style here] "
#    -> "Leadership Style to be applied: [Insert your leadership style here and why you chose it.] "
#    (appears twice: Key Issue 1 and Key Issue 2)
$d.Content.Find.Execute(
    "Leadership Style to be applied: [Insert your leadership style here] ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Leadership Style to be applied: [Insert your leadership style here and why you chose it.] ", 2
)

# 2) "Approach 1: How this approach motivates the team"
#    -> "Approach 1: Describe your approach and how this approach motivates the team"
#    (must run before the generic "How this approach motivates the team" replace below)
$d.Content.Find.Execute(
    "Approach 1: How this approach motivates the team", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Approach 1: Describe your approach and how this approach motivates the team", 2
)

# 3) Remaining bare "How this approach motivates the team"
#    -> "Describe your approach and how this approach motivates the team"
#    (appears 3 times: Key Issue 1 Approach 2, Key Issue 2 Approach 1, Key Issue 2 Approach 2)
$d.Content.Find.Execute(
    "How this approach motivates the team", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Describe your approach and how this approach motivates the team", 2
)
